# Auto-update draw results: append the 2025-11-26 Pick 3 draw as a new
# row (row 71) at the bottom of the Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow       = 71
$drawDate     = "2025-11-26"
$game         = "Pick 3"
$phase        = "251126"
$result       = "5-6-6"
$insertedAt   = "2025-11-26T21:40:21.376+04:00"

# Columns B (Game) and D (Result) don't look like numbers/dates, so a
# plain assignment keeps them as plain text.
$ws.Range("B$newRow").Value = $game
$ws.Range("D$newRow").Value = $result

# Columns A, C and E hold text that *looks* like a date / number / ISO
# timestamp ("2025-11-26", "251126", "2025-11-26T21:40:21.376+04:00").
# Excel's Value setter auto-converts such look-alikes into a real
# DateTime/Double. Prefixing with a single quote forces it to be kept
# as literal text (matching how the existing rows store these columns
# as strings), then ClearFormats() drops the transient "quote prefix"
# formatting so the new cells end up with the same (default) style as
# every other cell in the sheet.
$ws.Range("A$newRow").Value = "'" + $drawDate
$ws.Range("A$newRow").ClearFormats()

$ws.Range("C$newRow").Value = "'" + $phase
$ws.Range("C$newRow").ClearFormats()

$ws.Range("E$newRow").Value = "'" + $insertedAt
$ws.Range("E$newRow").ClearFormats()
